# Commit: "Update Excel SCD0011 until SCD0016"
#
# The test-case id for this scenario moves from the old "DGS-221" /
# "SCD0206" numbering scheme to the new "SCD0011-037" / "SCD0011" one:
#   - the worksheet is renamed SCD0206 -> SCD0011
#   - the TC_ID cell (B2) changes from "DGS-221" to "SCD0011-037"
#   - column B is widened so the longer id still fits
#   - the active selection moves from A2 to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab: SCD0206 -> SCD0011
$ws.Name = "SCD0011"

# TC_ID cell: DGS-221 -> SCD0011-037
$ws.Range("B2").Value = "SCD0011-037"

# Column B needs to be a bit wider to fit "SCD0011-037"
$ws.Columns.Item(2).ColumnWidth = 11.67

# Active cell/selection ends up on B3
$ws.Range("B3").Select()
